# New crime data collected - update the 110th Precinct CompStat weekly report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) -----------------------
# "Volume 30   Number  14" -> "...  15"
$volRng = $ws.Range("A8")
$volChars = $volRng.Characters(21, 2)
$volChars.Text = "15"

# "Report Covering the Week  4/3/2023  Through  4/9/2023"
#   -> "...4/10/2023  Through  4/16/2023"
$weekRng = $ws.Range("C9")
$startChars = $weekRng.Characters(27, 8)
$startChars.Text = "4/10/2023"
$endChars = $weekRng.Characters(47, 8)
$endChars.Text = "4/16/2023"

# --- Helper: style reference cells used to normalize formatting when a ----
# --- cell's type changes (number <-> text placeholder) --------------------
$numStyleSrc = $ws.Range("D15")   # plain integer style (s=15)
$pctStyleSrc = $ws.Range("E15")   # percentage style (s=16)
$txtStyleSrc = $ws.Range("D14")   # text placeholder style (s=14, shared string "0")

function Set-NumberCell($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-TextZeroCell($ws, $addr, $styleSrc) {
    $dst = $ws.Range($addr)
    $dst.Value = "'0"
    $styleSrc.Copy()
    $dst.PasteSpecial(-4122)
}

function Set-NumberCellWithStyle($ws, $addr, $value, $styleSrc) {
    $dst = $ws.Range($addr)
    $dst.Value = $value
    $styleSrc.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $value
}

# --- Row 15 (Murder) -------------------------------------------------------
Set-TextZeroCell $ws "C15" $txtStyleSrc
Set-NumberCell $ws "E15" -100
Set-NumberCell $ws "G15" 3
Set-NumberCell $ws "H15" 33.333333333333
Set-NumberCell $ws "J15" 8
Set-NumberCell $ws "K15" -12.5
Set-NumberCell $ws "L15" -50
Set-NumberCell $ws "M15" -30

# --- Row 16 (Rape) -----------------------------------------------------
Set-NumberCell $ws "C16" 7
Set-NumberCell $ws "D16" 4
Set-NumberCell $ws "E16" 75
Set-NumberCell $ws "F16" 27
Set-NumberCell $ws "H16" 12.5
Set-NumberCell $ws "I16" 104
Set-NumberCell $ws "J16" 86
Set-NumberCell $ws "K16" 20.930232558139
Set-NumberCell $ws "L16" 92.592592592592
Set-NumberCell $ws "M16" -13.333333333333
Set-NumberCell $ws "N16" -77.586206896551

# --- Row 17 (Robbery) -------------------------------------------------
Set-NumberCell $ws "C17" 17
Set-NumberCell $ws "D17" 16
Set-NumberCell $ws "E17" 6.25
Set-NumberCell $ws "F17" 49
Set-NumberCell $ws "G17" 35
Set-NumberCell $ws "H17" 40
Set-NumberCell $ws "I17" 180
Set-NumberCell $ws "J17" 138
Set-NumberCell $ws "K17" 30.434782608695
Set-NumberCell $ws "L17" 83.673469387755
Set-NumberCell $ws "M17" 140
Set-NumberCell $ws "N17" 53.846153846153

# --- Row 18 (Fel. Assault) ---------------------------------------------
Set-NumberCell $ws "C18" 2
Set-NumberCell $ws "D18" 2
Set-NumberCell $ws "E18" 0
Set-NumberCell $ws "F18" 14
Set-NumberCell $ws "G18" 10
Set-NumberCell $ws "H18" 40
Set-NumberCell $ws "I18" 61
Set-NumberCell $ws "J18" 42
Set-NumberCell $ws "K18" 45.238095238095
Set-NumberCell $ws "L18" 22
Set-NumberCell $ws "M18" -47.413793103448
Set-NumberCell $ws "N18" -91.562932226832

# --- Row 19 (Burglary) --------------------------------------------------
Set-NumberCell $ws "C19" 22
Set-NumberCell $ws "D19" 17
Set-NumberCell $ws "E19" 29.411764705882
Set-NumberCell $ws "F19" 76
Set-NumberCell $ws "G19" 60
Set-NumberCell $ws "H19" 26.666666666666
Set-NumberCell $ws "I19" 255
Set-NumberCell $ws "J19" 391
Set-NumberCell $ws "K19" -34.782608695652
Set-NumberCell $ws "L19" 114.285714285714
Set-NumberCell $ws "M19" 57.407407407407
Set-NumberCell $ws "N19" -19.558359621451

# --- Row 20 (Gr. Larceny) -----------------------------------------------
Set-NumberCell $ws "D20" 1
Set-NumberCell $ws "E20" 400
Set-NumberCell $ws "F20" 17
Set-NumberCell $ws "G20" 9
Set-NumberCell $ws "H20" 88.888888888888
Set-NumberCell $ws "I20" 80
Set-NumberCell $ws "J20" 47
Set-NumberCell $ws "K20" 70.212765957446
Set-NumberCell $ws "L20" 100
Set-NumberCell $ws "M20" 77.777777777777
Set-NumberCell $ws "N20" -87.117552334943

# --- Row 21 (G.L.A. / TOTAL, bold) --------------------------------------
Set-NumberCell $ws "C21" 53
Set-NumberCell $ws "D21" 41
Set-NumberCell $ws "E21" 29.268292682926
Set-NumberCell $ws "F21" 187
Set-NumberCell $ws "G21" 141
Set-NumberCell $ws "H21" 32.624113475177
Set-NumberCell $ws "I21" 689
Set-NumberCell $ws "J21" 713
Set-NumberCell $ws "K21" -3.366058906030
Set-NumberCell $ws "L21" 83.733333333333
Set-NumberCell $ws "M21" 30.245746691871
Set-NumberCell $ws "N21" -69.459219858156

# --- Row 22 (Transit) ---------------------------------------------------
Set-TextZeroCell $ws "C22" $txtStyleSrc
Set-NumberCellWithStyle $ws "D22" 1 $numStyleSrc
Set-NumberCellWithStyle $ws "E22" -100 $pctStyleSrc
Set-NumberCell $ws "F22" 3
Set-NumberCellWithStyle $ws "G22" 1 $numStyleSrc
Set-NumberCellWithStyle $ws "H22" 200 $pctStyleSrc
Set-NumberCell $ws "J22" 8
Set-NumberCell $ws "K22" 75
Set-NumberCell $ws "M22" 27.272727272727

# --- Row 24 (Petit Larceny) ---------------------------------------------
Set-NumberCell $ws "C24" 57
Set-NumberCell $ws "D24" 45
Set-NumberCell $ws "E24" 26.666666666666
Set-NumberCell $ws "F24" 202
Set-NumberCell $ws "G24" 184
Set-NumberCell $ws "H24" 9.782608695652
Set-NumberCell $ws "I24" 821
Set-NumberCell $ws "J24" 683
Set-NumberCell $ws "K24" 20.204978038067
Set-NumberCell $ws "L24" 82.039911308204
Set-NumberCell $ws "M24" 94.089834515366

# --- Row 25 (Misd. Assault) ---------------------------------------------
Set-NumberCell $ws "C25" 20
Set-NumberCell $ws "D25" 13
Set-NumberCell $ws "E25" 53.846153846153
Set-NumberCell $ws "F25" 89
Set-NumberCell $ws "G25" 69
Set-NumberCell $ws "H25" 28.985507246376
Set-NumberCell $ws "I25" 296
Set-NumberCell $ws "J25" 215
Set-NumberCell $ws "K25" 37.674418604651
Set-NumberCell $ws "L25" 51.020408163265
Set-NumberCell $ws "M25" 66.292134831460

# --- Row 26 (UCR Rape*) --------------------------------------------------
Set-TextZeroCell $ws "C26" $txtStyleSrc
Set-NumberCell $ws "E26" -100
Set-NumberCell $ws "G26" 4
Set-NumberCell $ws "H26" 0
Set-NumberCell $ws "J26" 12
Set-NumberCell $ws "K26" -16.666666666666
Set-NumberCell $ws "L26" -44.444444444444

# --- Row 27 (Other Sex Crimes) ------------------------------------------
Set-TextZeroCell $ws "C27" $txtStyleSrc
Set-NumberCell $ws "D27" 3
Set-NumberCell $ws "E27" -100
Set-NumberCell $ws "G27" 9
Set-NumberCell $ws "H27" 0
Set-NumberCell $ws "J27" 31
Set-NumberCell $ws "K27" 22.580645161290
Set-NumberCell $ws "L27" 31.034482758620
